# Tache.xlsx - "Correction + cahier des charges" commit
# Corrects spelling/grammar mistakes across the task list and fills in an
# "ok" remark for the orthography-correction task row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Column-legend row (row 2), edited first ---
$ws.Range("A2").Value = "Personne devant exécuter la tâche"
$ws.Range("B2").Value = "Description complète et précise de la tâche à effectuer"

# --- Header row (row 1) ---
$ws.Range("B1").Value = "Description de tâche"
$ws.Range("D1").Value = "Tâche donnée par "

# --- Rest of column-legend row (row 2) ---
$ws.Range("D2").Value = "Personne donnant la tâche à exécuter."
$ws.Range("E2").Value = "Noter : OK ! Si aucun problème n'est rencontré."

# --- Row 3 ---
$ws.Range("B3").Value = "Rédaction d'un Background Complet donnant une ligne directive au projet. Quelle est la place du joueur dans le jeu, ce qu'il représente … Quel est le but du personnage, qui il est, et comment arrive t'il à ses fins"
$ws.Range("E3").Value = "Problème de compréhension dans la tache : Réponses pas assez précises."

# --- Row 4 ---
$ws.Range("B4").Value = "Donner les premières mécaniques de jeu, comment s'enchaine un ""Monde"" , comment le joueur intéragit avec le personnage. "

# --- Row 5 ---
$ws.Range("B5").Value = "Création d'un premier ""Level"". Mécanique complète, design complet, enchaînement des évènements complets."

# --- Row 6 ---
$ws.Range("B6").Value = "Création d'un écran titre. Personnage complet en action avec un décord. Le tous en 1024 * 768"

# --- Row 8 ---
$ws.Range("B8").Value = "Création d'un Fichier de contact regroupant les informations personnelles de chacun des membre du groupe"

# --- Row 9 ---
$ws.Range("B9").Value = "Composition / Recherche d'une musique décrivant l'univers du jeu"

# --- Row 10 ---
$ws.Range("B10").Value = "Création d'un premier moteur physique dans le jeu"

# --- Row 11 ---
$ws.Range("B11").Value = "Mise en place d'une architecture de travail au sein du groupe de développement"

# --- Row 12 ---
$ws.Range("B12").Value = "Création de l'emploi du temps de travail des développeurs."

# --- Row 15 ---
$ws.Range("B15").Value = "Création d'une interface de jeu + Menu principal"

# --- Row 16 ---
$ws.Range("B16").Value = "Liste des différentes actions possibles par le joueur / personnage"

# --- Row 17 ---
$ws.Range("B17").Value = "Création de 10 mécanismes (Pièges) différents"

# --- Row 19 ---
$ws.Range("B19").Value = "Création de la première musique d'ambience : menu principal"

# --- Row 20 : cahier des charges correction note ---
$ws.Range("B20").Value = "Correction orthographe / syntaxe de la totalité du Repo"

# --- Row 24 : mark the orthography task as ok ---
$ws.Range("E24").Value = "ok"

# --- Restore the view (scroll position + selection) ---
$ws.Application.ActiveWindow.ScrollRow = 23
$sheetView = $ws.Application.ActiveWindow
$ws.Range("A23").Select()
$ws.Range("E25").Select()
